$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leaderboard")

$ws.Range("B2").Value = 605
$ws.Range("C2").Value = 464
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 7
